$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.000.86'
$ws.Range("E2").Value = '  -2.02%  '
$ws.Range("D3").Value = '2.677.41'
$ws.Range("E3").Value = '  -3.37%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '546.79'
$ws.Range("E5").Value = '  -4.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.53'
$ws.Range("E6").Value = '  -2.13%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("E9").Value = '  -5.08%  '
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.10'
$ws.Range("E12").Value = '  -13.30%  '
$ws.Range("D13").Value = '3.150.22'
$ws.Range("E13").Value = '  -3.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.07'
$ws.Range("E14").Value = '  -4.72%  '
$ws.Range("D15").Value = '62.874.75'
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("E16").Value = '  -5.37%  '
$ws.Range("D17").Value = '2.677.32'
$ws.Range("E17").Value = '  -3.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.91'
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("E19").Value = '  -6.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.92'
$ws.Range("E20").Value = '  -5.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.29'
$ws.Range("E21").Value = '  -6.20%  '
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").Value = '  -5.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.30'
$ws.Range("E24").Value = '  -2.92%  '
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.11'
$ws.Range("E27").Value = '  -6.17%  '
$ws.Range("D28").Value = '0.0₃0852'
$ws.Range("E28").Value = '  -7.28%  '
$ws.Range("E29").Value = '  -3.77%  '
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.99'
$ws.Range("E31").Value = '  -5.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '166.02'
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.80'
$ws.Range("E34").Value = '  -4.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.49'
$ws.Range("E35").Value = '  -3.76%  '
$ws.Range("E36").Value = '  -6.67%  '
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '339.37'
$ws.Range("E38").Value = '  -3.19%  '
$ws.Range("B39").Value = 'SuiNetwork'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.934'
$ws.Range("E39").Value = '  -7.68%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.13'
$ws.Range("E40").Value = '  -3.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.02'
$ws.Range("E41").Value = '  -2.96%  '
$ws.Range("E42").Value = '  -7.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.29'
$ws.Range("E43").Value = '  -6.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.69'
$ws.Range("E44").Value = '  -8.96%  '
$ws.Range("E45").Value = '  -6.29%  '
$ws.Range("E46").Value = '  -3.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.998'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("E49").Value = '  -4.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.75'
$ws.Range("E50").Value = '  -6.29%  '
$ws.Range("E51").Value = '  -6.16%  '
